$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.026551335456299
$ws.Range("D2").Value = 1.037378488097944
$ws.Range("E2").Value = 1.047964084582233
$ws.Range("F2").Value = 1.051988769712932
$ws.Range("I2").Value = 1.03597085415986
$ws.Range("J2").Value = 1.031714145944161
$ws.Range("K2").Value = 1.040169542702951
$ws.Range("L2").Value = 1.050725252719756
$ws.Range("M2").Value = 1.054738745987476
$ws.Range("N2").Value = 1.01460960757403
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.027367332931989
$ws.Range("D3").Value = 1.03800318373375
$ws.Range("E3").Value = 1.048799478953313
$ws.Range("F3").Value = 1.052782148601883
$ws.Range("I3").Value = 1.03611461768672
$ws.Range("J3").Value = 1.032170638673613
$ws.Range("K3").Value = 1.04060452395568
$ws.Range("L3").Value = 1.051372536393778
$ws.Range("M3").Value = 1.055344930777183
$ws.Range("N3").Value = 1.014760818067133
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.027895795192256
$ws.Range("D4").Value = 1.038407506053869
$ws.Range("E4").Value = 1.049341137359727
$ws.Range("F4").Value = 1.053296156182433
$ws.Range("I4").Value = 1.036206033259433
$ws.Range("J4").Value = 1.032465820402144
$ws.Range("K4").Value = 1.04088539762326
$ws.Range("L4").Value = 1.051791825357059
$ws.Range("M4").Value = 1.055737153455715
$ws.Range("N4").Value = 1.014858577831136
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.028118068554504
$ws.Range("D5").Value = 1.038577505802278
$ws.Range("E5").Value = 1.049569112481967
$ws.Range("F5").Value = 1.053512395768915
$ws.Range("I5").Value = 1.036244078248415
$ws.Range("J5").Value = 1.032589865743561
$ws.Range("K5").Value = 1.041003334648157
$ws.Range("L5").Value = 1.051968201256345
$ws.Range("M5").Value = 1.055902037518519
$ws.Range("N5").Value = 1.014899655588891
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.02815539550497
$ws.Range("D6").Value = 1.038606050781004
$ws.Range("E6").Value = 1.049607405809939
$ws.Range("F6").Value = 1.053548712130793
$ws.Range("I6").Value = 1.036250443500832
$ws.Range("J6").Value = 1.03261069058836
$ws.Range("K6").Value = 1.041023128409518
$ws.Range("L6").Value = 1.051997821761404
$ws.Range("M6").Value = 1.05592972185615
$ws.Range("N6").Value = 1.014906551513089
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.02789876479819
$ws.Range("D7").Value = 1.03840977751223
$ws.Range("E7").Value = 1.04934418254777
$ws.Range("F7").Value = 1.053299044996393
$ws.Range("I7").Value = 1.036206543136761
$ws.Range("J7").Value = 1.032467478096434
$ws.Range("K7").Value = 1.040886974064363
$ws.Range("L7").Value = 1.051794181682828
$ws.Range("M7").Value = 1.055739356671735
$ws.Range("N7").Value = 1.014859126795053
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.026827009934255
$ws.Range("D8").Value = 1.037589584562851
$ws.Range("E8").Value = 1.048246180571802
$ws.Range("F8").Value = 1.052256762420398
$ws.Range("I8").Value = 1.036019772333623
$ws.Range("J8").Value = 1.031868460360877
$ws.Range("K8").Value = 1.040316667606224
$ws.Range("L8").Value = 1.050943910441769
$ws.Range("M8").Value = 1.054943611972338
$ws.Range("N8").Value = 1.014660726895675
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.024942016726123
$ws.Range("D9").Value = 1.036145167771577
$ws.Range("E9").Value = 1.046319891457189
$ws.Range("F9").Value = 1.050425106293291
$ws.Range("I9").Value = 1.035678375028723
$ws.Range("J9").Value = 1.030811442139854
$ws.Range("K9").Value = 1.039307272066053
$ws.Range("L9").Value = 1.049449167890069
$ws.Range("M9").Value = 1.053541334584062
$ws.Range("N9").Value = 1.014310500902465
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.023687856997155
$ws.Range("D10").Value = 1.035182918789916
$ws.Range("E10").Value = 1.045041552418459
$ws.Range("F10").Value = 1.049207470487376
$ws.Range("I10").Value = 1.035442571982564
$ws.Range("J10").Value = 1.030105848885493
$ws.Range("K10").Value = 1.038631439241546
$ws.Range("L10").Value = 1.048455154608879
$ws.Range("M10").Value = 1.052606534644644
$ws.Range("N10").Value = 1.014076626620213
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.023145406108591
$ws.Range("D11").Value = 1.034766442908009
$ws.Range("E11").Value = 1.044489430848297
$ws.Range("F11").Value = 1.048681068797431
$ws.Range("I11").Value = 1.035338531821085
$ws.Range("J11").Value = 1.02980011818577
$ws.Range("K11").Value = 1.038338125496905
$ws.Range("L11").Value = 1.048025345373421
$ws.Range("M11").Value = 1.052201789472839
$ws.Range("N11").Value = 1.013975269488763
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.022944008767193
$ws.Range("D12").Value = 1.034611774806706
$ws.Range("E12").Value = 1.044284561819893
$ws.Range("F12").Value = 1.048485668422412
$ws.Range("I12").Value = 1.035299596706267
$ws.Range("J12").Value = 1.029686526785992
$ws.Range("K12").Value = 1.038229075917796
$ws.Range("L12").Value = 1.047865787673861
$ws.Range("M12").Value = 1.052051455215401
$ws.Range("N12").Value = 1.013937608159709
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.022987204942527
$ws.Range("D13").Value = 1.034644950270264
$ws.Range("E13").Value = 1.044328497227942
$ws.Range("F13").Value = 1.048527576611539
$ws.Range("I13").Value = 1.035307961528077
$ws.Range("J13").Value = 1.029710893837156
$ws.Range("K13").Value = 1.038252471915898
$ws.Range("L13").Value = 1.04790000912528
$ws.Range("M13").Value = 1.052083702119586
$ws.Range("N13").Value = 1.013945687213678
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.023128756633395
$ws.Range("D14").Value = 1.034753657380024
$ws.Range("E14").Value = 1.044472491943375
$ws.Range("F14").Value = 1.04866491430519
$ws.Range("I14").Value = 1.035335319342503
$ws.Range("J14").Value = 1.029790729277763
$ws.Range("K14").Value = 1.038329113452276
$ws.Range("L14").Value = 1.048012154385721
$ws.Range("M14").Value = 1.052189362661498
$ws.Range("N14").Value = 1.01397215664963
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.023215983683936
$ws.Range("D15").Value = 1.034820639395112
$ws.Range("E15").Value = 1.044561240191945
$ws.Range("F15").Value = 1.048749549699733
$ws.Range("I15").Value = 1.035352136997763
$ws.Range("J15").Value = 1.029839914669029
$ws.Range("K15").Value = 1.038376321650897
$ws.Range("L15").Value = 1.048081263100979
$ws.Range("M15").Value = 1.052254464470679
$ws.Range("N15").Value = 1.013988463662932
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.023723870637848
$ws.Range("D16").Value = 1.03521056296734
$ws.Range("E16").Value = 1.045078224739137
$ws.Range("F16").Value = 1.049242423962793
$ws.Range("I16").Value = 1.035449436073584
$ws.Range("J16").Value = 1.030126135031935
$ws.Range("K16").Value = 1.038650891436331
$ws.Range("L16").Value = 1.048483692524173
$ws.Range("M16").Value = 1.052633397013389
$ws.Range("N16").Value = 1.014083351542842
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.024042618792227
$ws.Range("D17").Value = 1.035455202401518
$ws.Range("E17").Value = 1.045402893786771
$ws.Range("F17").Value = 1.049551818017622
$ws.Range("I17").Value = 1.035509951542506
$ws.Range("J17").Value = 1.030305619739249
$ws.Range("K17").Value = 1.038822942449257
$ws.Range("L17").Value = 1.04873628905338
$ws.Range("M17").Value = 1.052871100456495
$ws.Range("N17").Value = 1.014142848952486
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.024228597803963
$ws.Range("D18").Value = 1.035597914097544
$ws.Range("E18").Value = 1.045592403433707
$ws.Range("F18").Value = 1.049732363501317
$ws.Range("I18").Value = 1.03554506234794
$ws.Range("J18").Value = 1.030410290348593
$ws.Range("K18").Value = 1.038923231799863
$ws.Range("L18").Value = 1.048883682603567
$ws.Range("M18").Value = 1.053009751572745
$ws.Range("N18").Value = 1.014177544261966
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.024292021731443
$ws.Range("D19").Value = 1.035646578034274
$ws.Range("E19").Value = 1.045657044278888
$ws.Range("F19").Value = 1.049793938550216
$ws.Range("I19").Value = 1.035557002506541
$ws.Range("J19").Value = 1.030445976927152
$ws.Range("K19").Value = 1.038957416792407
$ws.Range("L19").Value = 1.048933949827287
$ws.Range("M19").Value = 1.053057028429577
$ws.Range("N19").Value = 1.014189373006632
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.024008414058088
$ws.Range("D20").Value = 1.035428953065817
$ws.Range("E20").Value = 1.045368045821541
$ws.Range("F20").Value = 1.049518614553409
$ws.Range("I20").Value = 1.035503478125174
$ws.Range("J20").Value = 1.030286364765947
$ws.Range("K20").Value = 1.038804489726618
$ws.Range("L20").Value = 1.048709181797235
$ws.Range("M20").Value = 1.052845596833021
$ws.Range("N20").Value = 1.014136466321847
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.023087070606098
$ws.Range("D21").Value = 1.034721645013746
$ws.Range("E21").Value = 1.044430083178641
$ws.Range("F21").Value = 1.048624468209615
$ws.Range("I21").Value = 1.035327271145481
$ws.Range("J21").Value = 1.029767220531908
$ws.Range("K21").Value = 1.038306547173511
$ws.Range("L21").Value = 1.047979127832299
$ws.Range("M21").Value = 1.052158248082766
$ws.Range("N21").Value = 1.013964362409613
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.022508324704895
$ws.Range("D22").Value = 1.03427710447273
$ws.Range("E22").Value = 1.043841585114533
$ws.Range("F22").Value = 1.048063028260474
$ws.Range("I22").Value = 1.035214805600365
$ws.Range("J22").Value = 1.029440644443335
$ws.Range("K22").Value = 1.037992894245075
$ws.Range("L22").Value = 1.047520650443492
$ws.Range("M22").Value = 1.051726120897681
$ws.Range("N22").Value = 1.013856080142989
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.022815077189635
$ws.Range("D23").Value = 1.034512746878462
$ws.Range("E23").Value = 1.044153441214744
$ws.Range("F23").Value = 1.048360586866034
$ws.Range("I23").Value = 1.035274584402418
$ws.Range("J23").Value = 1.029613784251177
$ws.Range("K23").Value = 1.03815922178765
$ws.Range("L23").Value = 1.047763646565667
$ws.Range("M23").Value = 1.051955195688701
$ws.Range("N23").Value = 1.0139134894588
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.024023869524225
$ws.Range("D24").Value = 1.03544081395642
$ws.Range("E24").Value = 1.045383791698026
$ws.Range("F24").Value = 1.049533617519123
$ws.Range("I24").Value = 1.035506403761234
$ws.Range("J24").Value = 1.03029506532084
$ws.Range("K24").Value = 1.038812827918987
$ws.Range("L24").Value = 1.048721430219636
$ws.Range("M24").Value = 1.052857120814481
$ws.Range("N24").Value = 1.014139350384283
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.025428898778397
$ws.Range("D25").Value = 1.036518469714025
$ws.Range("E25").Value = 1.04681686025155
$ws.Range("F25").Value = 1.050898030821224
$ws.Range("I25").Value = 1.035768084499712
$ws.Range("J25").Value = 1.031084873383127
$ws.Range("K25").Value = 1.039568742762129
$ws.Range("L25").Value = 1.049835164477034
$ws.Range("M25").Value = 1.053903855185598
$ws.Range("N25").Value = 1.014401113449081
